# Update the probability-threshold table values (columns B/C, rows 2-3)
# on Sheet1, then move the cell selection to C3 to match the saved
# workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.8
$ws.Range("C2").Value = 10.7
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 9.3000000000000007

$ws.Range("C3").Select()
